$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2999.4
$ws.Range("J40").Value = 2999.25
$ws.Range("L40").Value = 2999.25
$ws.Range("N40").Value = -3349.25

$ws.Range("H112").Value = 1594.3729
$ws.Range("I112").Value = 800
$ws.Range("J112").Value = 1608.069
$ws.Range("K112").Value = 2400
$ws.Range("L112").Value = 4824.207
$ws.Range("M112").Value = -1292
$ws.Range("N112").Value = -7040.207

$ws.Range("H113").Value = 10705.333
$ws.Range("I113").Value = 14886.25
$ws.Range("J113").Value = 2343.5
$ws.Range("K113").Value = 14886.25
$ws.Range("L113").Value = 2343.5
$ws.Range("M113").Value = -11632.25
$ws.Range("N113").Value = -8851.5

$ws.Range("H129").Value = 1084
$ws.Range("I129").Value = 765.25
$ws.Range("J129").Value = 1113.6511
$ws.Range("K129").Value = 2295.75
$ws.Range("L129").Value = 3340.9533
$ws.Range("M129").Value = 2704.25
$ws.Range("N129").Value = -13340.9533

$ws.Range("H137").Value = 68227.8
$ws.Range("I137").Value = 866.6667
$ws.Range("J137").Value = 113135.22
$ws.Range("K137").Value = 2600.0001
$ws.Range("L137").Value = 339405.66
$ws.Range("M137").Value = -50.0001000000002
$ws.Range("N137").Value = -344505.66

$ws.Range("H138").Value = 1780.102
$ws.Range("I138").Value = 1307.4138
$ws.Range("J138").Value = 1978.7681
$ws.Range("K138").Value = 3922.2414
$ws.Range("L138").Value = 5936.3043
$ws.Range("M138").Value = 1217.7586
$ws.Range("N138").Value = -16216.3043

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38915.637
$ws.Range("I61").Value = 48181.293
$ws.Range("J61").Value = 7412.4
$ws.Range("K61").Value = 48181.293
$ws.Range("L61").Value = 7412.4
$ws.Range("M61").Value = -47969.293
$ws.Range("N61").Value = -7836.4

$ws.Range("H74").Value = 767.70966
$ws.Range("I74").Value = 767.70966
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 767.70966
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 106.29034
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 767.70966
$ws.Range("I77").Value = 767.70966
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3838.5483
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 529.4517000000001
$ws.Range("N77").ClearContents()

$ws.Range("H88").Value = 2835.5833
$ws.Range("I88").Value = 2124.5
$ws.Range("K88").Value = 2124.5
$ws.Range("M88").Value = -1718.5

$ws.Range("H91").Value = 2835.5833
$ws.Range("I91").Value = 2124.5
$ws.Range("K91").Value = 2124.5
$ws.Range("M91").Value = -720.5

$ws.Range("H132").Value = 1598.0441
$ws.Range("I132").Value = 1075.3572
$ws.Range("J132").Value = 2442.3845
$ws.Range("K132").Value = 3226.0716
$ws.Range("L132").Value = 7327.1535
$ws.Range("M132").Value = -696.0715999999998
$ws.Range("N132").Value = -12387.1535

$ws.Range("H136").Value = 38915.637
$ws.Range("I136").Value = 48181.293
$ws.Range("J136").Value = 7412.4
$ws.Range("K136").Value = 144543.879
$ws.Range("L136").Value = 22237.2
$ws.Range("M136").Value = -141993.879
$ws.Range("N136").Value = -27337.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4112.857
$ws.Range("I20").Value = 2932.3333
$ws.Range("K20").Value = 2932.3333
$ws.Range("M20").Value = -2685.3333

$ws.Range("H86").Value = 801379.6
$ws.Range("I86").Value = 1949
$ws.Range("J86").Value = 1334333.4
$ws.Range("K86").Value = 1949
$ws.Range("L86").Value = 1334333.4
$ws.Range("M86").Value = -826
$ws.Range("N86").Value = -1336579.4

$ws.Range("H89").Value = 801379.6
$ws.Range("I89").Value = 1949
$ws.Range("J89").Value = 1334333.4
$ws.Range("K89").Value = 9745
$ws.Range("L89").Value = 6671667
$ws.Range("M89").Value = -4129
$ws.Range("N89").Value = -6682899

$ws.Range("H94").Value = 345.54544
$ws.Range("I94").Value = 332.55554
$ws.Range("J94").Value = 404
$ws.Range("K94").Value = 332.55554
$ws.Range("L94").Value = 404
$ws.Range("M94").Value = 118.44446
$ws.Range("N94").Value = -1306

$ws.Range("H99").Value = 1775.2
$ws.Range("I99").Value = 1718.4166
$ws.Range("J99").Value = 2002.3334
$ws.Range("K99").Value = 1718.4166
$ws.Range("L99").Value = 2002.3334
$ws.Range("M99").Value = -220.4166
$ws.Range("N99").Value = -4998.3334

$ws.Range("H105").Value = 2730.182
$ws.Range("I105").Value = 2526.8572
$ws.Range("K105").Value = 2526.8572
$ws.Range("M105").Value = -779.8571999999999

$ws.Range("H107").Value = 1790.6923
$ws.Range("I107").Value = 1688.8334
$ws.Range("K107").Value = 1688.8334
$ws.Range("M107").Value = 231.1666

$ws.Range("H134").Value = 3973.561
$ws.Range("I134").Value = 3748.5
$ws.Range("J134").Value = 4773.778
$ws.Range("K134").Value = 11245.5
$ws.Range("L134").Value = 14321.334
$ws.Range("M134").Value = -8710.5
$ws.Range("N134").Value = -19391.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3029.5334
$ws.Range("I31").Value = 2145.5557
$ws.Range("K31").Value = 2145.5557
$ws.Range("M31").Value = -1850.5557

$ws.Range("H34").Value = 3029.5334
$ws.Range("I34").Value = 2145.5557
$ws.Range("K34").Value = 2145.5557
$ws.Range("M34").Value = -1943.5557

$ws.Range("H92").Value = 44995
$ws.Range("J92").Value = 44995
$ws.Range("L92").Value = 44995
$ws.Range("N92").Value = -49987

$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 126.5
$ws.Range("J12").Value = 139.5
$ws.Range("L12").Value = 418.5
$ws.Range("N12").Value = -764.5

$ws.Range("H33").Value = 196.28572
$ws.Range("I33").Value = 125.57143
$ws.Range("K33").Value = 753.42858
$ws.Range("M33").Value = -470.42858

$ws.Range("H62").Value = 4062.5
$ws.Range("I62").Value = 4250
$ws.Range("K62").Value = 12750
$ws.Range("M62").Value = -12064

$ws.Range("H65").Value = 4062.5
$ws.Range("I65").Value = 4250
$ws.Range("K65").Value = 38250
$ws.Range("M65").Value = -34818

$ws.Range("H131").Value = 17984.277
$ws.Range("I131").Value = 418.2
$ws.Range("J131").Value = 20075.477
$ws.Range("K131").Value = 1254.6
$ws.Range("L131").Value = 60226.431
$ws.Range("M131").Value = 3785.4
$ws.Range("N131").Value = -70306.431

$ws.Range("H137").Value = 3776.7368
$ws.Range("I137").Value = 2553.7144
$ws.Range("J137").Value = 4052.9033
$ws.Range("K137").Value = 7661.1432
$ws.Range("L137").Value = 12158.7099
$ws.Range("M137").Value = -2561.1432
$ws.Range("N137").Value = -22358.7099

$ws.Range("H141").Value = 2939.5
$ws.Range("J141").Value = 2360.8333
$ws.Range("L141").Value = 7082.499899999999
$ws.Range("N141").Value = -17442.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H92").Value = 23750
$ws.Range("J92").Value = 23750
$ws.Range("L92").Value = 23750
$ws.Range("N92").Value = -27494

$ws.Range("H97").Value = 1121.9231
$ws.Range("I97").Value = 1132
$ws.Range("J97").Value = 1099.25
$ws.Range("K97").Value = 1132
$ws.Range("L97").Value = 1099.25
$ws.Range("M97").Value = -636
$ws.Range("N97").Value = -2091.25

$ws.Range("H102").Value = 9999
$ws.Range("I102").Value = 9999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 9999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -8377
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 2052.3333
$ws.Range("J122").Value = 3157
$ws.Range("L122").Value = 9471
$ws.Range("N122").Value = -14371

$ws.Range("H132").Value = 858867.1
$ws.Range("I132").Value = 1016275.3
$ws.Range("K132").Value = 3048825.9
$ws.Range("M132").Value = -3046295.9

$ws.Range("H135").Value = 63240
$ws.Range("J135").Value = 63240
$ws.Range("L135").Value = 63240
$ws.Range("N135").Value = -73380

$ws.Range("H136").Value = 19020
$ws.Range("J136").Value = 19020
$ws.Range("L136").Value = 57060
$ws.Range("N136").Value = -62160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9886.478999999999
$ws.Range("I40").Value = 9344.5
$ws.Range("K40").Value = 9344.5
$ws.Range("M40").Value = -9208.5
